# "Expanded allocate to move across all rows and columns"
#
# 1. Fill in the "min" sheet's allocation numbers (cols C:I) for every row,
#    and add the K column row-total formula.
# 2. Add a new "diff" sheet (current - min) for every row/column.
# 3. Update the remembered cell selections on "current" and "min".
# 4. Give column J (emer_drivers) a comfortable width where it is now
#    visible in the grid.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("current")
$ws2 = $wb.Worksheets.Item("min")

# ---------------------------------------------------------------------------
# 1. "min" sheet allocation numbers: (row, col, value) for columns C(3)..I(9)
# ---------------------------------------------------------------------------
$minData = @(
    @(2,3,0), @(2,4,0), @(2,5,5), @(2,6,0), @(2,7,1), @(2,8,5),
    @(2,9,1), @(3,6,7), @(3,9,2), @(4,3,72), @(4,5,148), @(4,6,4),
    @(4,7,0), @(4,8,0), @(4,9,16), @(5,3,0), @(5,4,6), @(5,5,4),
    @(5,6,0), @(5,7,1), @(5,8,1), @(5,9,1), @(6,6,2), @(7,3,6),
    @(7,5,6), @(7,6,1), @(7,9,1), @(8,3,4), @(8,5,8), @(8,6,1),
    @(8,9,1), @(9,3,2), @(11,3,5), @(11,5,10), @(11,6,1), @(11,9,1),
    @(12,8,2), @(13,3,1), @(13,5,2), @(13,9,1)
)

foreach ($item in $minData) {
    $ws2.Cells.Item($item[0], $item[1]).Value = $item[2]
}

# K column: total allocated per row
for ($r = 2; $r -le 13; $r++) {
    $ws2.Cells.Item($r, 11).Formula = "=SUM(C$r`:I$r)"
}

# Column J (emer_drivers) is now wide enough to read comfortably.
$ws2.Columns.Item(10).ColumnWidth = 12

# ---------------------------------------------------------------------------
# 2. New "diff" sheet = current - min, for every row/column
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "diff"

$headers = @("svc","priority","cdl","drivers","laborers","admin","cashier","equip_op","supers","emer_drivers")
for ($c = 1; $c -le $headers.Length; $c++) {
    $ws3.Cells.Item(1, $c).Value = $headers[$c - 1]
}

$rowInfo = @(
    @(2,"qrl",1), @(3,"payroll",1), @(4,"routine",2), @(5,"nwts",2),
    @(6,"office",3), @(7,"conv_ctr",3), @(8,"ds_da",3), @(9,"mss",3),
    @(10,"prop",4), @(11,"spec",4), @(12,"marine",4), @(13,"night",4)
)

$colLetters = @("A","B","C","D","E","F","G","H","I","J")

foreach ($info in $rowInfo) {
    $r = $info[0]
    $ws3.Cells.Item($r, 1).Value = $info[1]
    $ws3.Cells.Item($r, 2).Value = $info[2]
    for ($c = 3; $c -le 10; $c++) {
        $letter = $colLetters[$c - 1]
        $ws3.Cells.Item($r, $c).Formula = "=current!" + $letter + "$r-min!" + $letter + "$r"
    }
}

$ws3.Columns.Item(10).ColumnWidth = 12

# ---------------------------------------------------------------------------
# 3. Selections remembered per-sheet. Select on the background sheets first,
#    then re-activate "current" last so it stays the visible tab.
# ---------------------------------------------------------------------------
$ws3.Range("I18").Select() | Out-Null
$ws2.Range("C4").Select() | Out-Null

$ws1.Activate()
$ws1.Range("L7").Select() | Out-Null
